$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.594.48"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "3.713.38"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "622.65"
$ws.Range("E5").Value = "  +8.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "195.00"
$ws.Range("E6").Value = "  +13.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.636"
$ws.Range("E7").Value = "  +2.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.995"
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.731"
$ws.Range("E9").Value = "  +4.10%  "
$ws.Range("B10").Value = "Avalanche"
$ws.Range("C10").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "60.89"
$ws.Range("E10").Value = "  +18.95%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.161"
$ws.Range("E11").Value = "  -0.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000289"
$ws.Range("E12").Value = "  -1.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.49"
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("D14").Value = "4.299.69"
$ws.Range("E14").Value = "  +1.29%  "
$ws.Range("D15").Value = "3.709.21"
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.50"
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("E17").Value = "  +3.04%  "
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.91"
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("D20").Value = "68.477.49"
$ws.Range("E20").Value = "  +1.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "410.41"
$ws.Range("E21").Value = "  +0.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.70"
$ws.Range("E22").Value = "  +5.92%  "
$ws.Range("E23").Value = "  +3.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.10"
$ws.Range("E24").Value = "  +2.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.61"
$ws.Range("E25").Value = "  +8.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.12"
$ws.Range("E26").Value = "  +3.06%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.80"
$ws.Range("E27").Value = "  +2.19%  "
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.03"
$ws.Range("E28").Value = "  +0.89%  "
$ws.Range("E29").Value = "  +2.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.89"
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.66"
$ws.Range("E31").Value = "  +1.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.77"
$ws.Range("E32").Value = "  +2.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "46.68"
$ws.Range("E33").Value = "  +8.08%  "
$ws.Range("E34").Value = "  +5.95%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "68.17"
$ws.Range("E35").Value = "  +4.82%  "
$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "636.59"
$ws.Range("E36").Value = "  +5.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.415"
$ws.Range("E37").Value = "  +4.77%  "
$ws.Range("D38").Value = "0.0₃0829"
$ws.Range("E38").Value = "  -7.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("E40").Value = "  +0.28%  "
$ws.Range("E41").Value = "  +5.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.05"
$ws.Range("E42").Value = "  +2.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0447"
$ws.Range("E43").Value = "  +2.77%  "
$ws.Range("E44").Value = "  -1.25%  "
$ws.Range("D45").Value = "2.927.31"
$ws.Range("E45").Value = "  +6.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.46"
$ws.Range("E46").Value = "  +3.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.140"
$ws.Range("E47").Value = "  +4.77%  "
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "146.75"
$ws.Range("E49").Value = "  +2.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.06"
$ws.Range("E50").Value = "  -3.03%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.55"
$ws.Range("E51").Value = "  -13.58%  "
